# Apply the "Added reporting section Credentials in config file" change:
#  - Update the "admin" sheet (sheet1) test data (PAN / company names)
#  - Leave "login" sheet (sheet2) content as-is
#  - Add three new sheets: "business", "sub", "datalist" with new test data
#  - Adjust selections / active sheet to match the authored workbook

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1 "admin" - replace company pan / name test rows
# ---------------------------------------------------------------------
$admin = $wb.Worksheets.Item(1)

$admin.Range("A1").Value = "Company PAN"
$admin.Range("B1").Value = "Company Name"

$admin.Range("A2").Value = "CMNPK4241N"
$admin.Range("B2").Value = "company1"

$admin.Range("A3").Value = "KMNPK8241P"
$admin.Range("B3").Value = "company2"

$admin.Range("A4").Value = "KMNPK9241R"
$admin.Range("B4").Value = "company3"

$admin.Range("A5").Value = "KRPPK9241N"
$admin.Range("B5").Value = "company4"

$admin.Activate()
$admin.Range("B9").Select()

# ---------------------------------------------------------------------
# Sheet2 "login" - content unchanged, only the selection/active state
# moves away from it (handled automatically once another sheet becomes
# active below).
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Sheet3 "business" - new sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$business = $wb.Worksheets.Add($null, $lastSheet)
$business.Name = "business"

$admin.Range("A1:B1").Copy()
$business.Range("A1:C1").PasteSpecial(-4122)
$business.Application.CutCopyMode = $false

$business.Range("A1").Value = "Company Name"
$business.Range("B1").Value = "Business Unit"

$business.Range("A2").Value = "company1"
$business.Range("B2").Value = "unit01"

$business.Columns.Item(1).ColumnWidth = 16.0833333333335
$business.Columns.Item(2).ColumnWidth = 25.75
$business.Columns.Item(3).ColumnWidth = 24.5833333333337

# ---------------------------------------------------------------------
# Sheet4 "sub" - new sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sub = $wb.Worksheets.Add($null, $lastSheet)
$sub.Name = "sub"

$admin.Range("A1:B1").Copy()
$sub.Range("A1:C1").PasteSpecial(-4122)
$sub.Application.CutCopyMode = $false

$sub.Range("A1").Value = "Company Name"
$sub.Range("B1").Value = "Business Unit Name"
$sub.Range("C1").Value = "Sub Business Unit Name"

$sub.Range("A2").Value = "company1"
$sub.Range("B2").Value = "Unit01"
$sub.Range("C2").Value = "sub01"

$sub.Columns.Item(1).ColumnWidth = 16.7499999999995
$sub.Columns.Item(2).ColumnWidth = 23.0833333333334
$sub.Columns.Item(3).ColumnWidth = 24.7499999999997

$sub.Activate()
$sub.Range("A3").Select()

# ---------------------------------------------------------------------
# Sheet5 "datalist" - new sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$datalist = $wb.Worksheets.Add($null, $lastSheet)
$datalist.Name = "datalist"

$admin.Range("A1:B1").Copy()
$datalist.Range("A1:B1").PasteSpecial(-4122)
$datalist.Application.CutCopyMode = $false

$datalist.Range("A1").Value = "Source"
$datalist.Range("B1").Value = "Category"

$datalist.Range("A2").Value = "SAP"

$datalist.Activate()
$datalist.Range("A2").Select()
